$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 0.01
$ws.Range("B6").Value = -263724.5473518897
$ws.Range("B7").Value = 11018349.94395913
$ws.Range("B8").Value = 28330649.05696066
$ws.Range("B10").Value = 1336590.354541983

# --- Fed-in Capacity sheet ---
$ws = $wb.Worksheets.Item("Fed-in Capacity")
$ws.Range("J2").Value = 169.0966151720738
$ws.Range("R2").Value = 65.71641987298243
$ws.Range("I3").Value = 10.12574714858493
$ws.Range("J3").Value = 126.0910353404088
$ws.Range("K3").Value = 137.841438974359
$ws.Range("Q3").Value = 139.9817740860215
$ws.Range("R3").Value = 45.52166981132082
$ws.Range("J4").Value = 33.63624132272333
$ws.Range("K4").Value = 106.7437663446525
$ws.Range("P4").Value = 135.0065633140411
$ws.Range("Q4").Value = 65.34295837775146
$ws.Range("J5").Value = 169.0966151720738
$ws.Range("L5").Value = 235.7664149699872
$ws.Range("O5").Value = 230.0982114216867
$ws.Range("P5").Value = 231.2329957552695
$ws.Range("Q5").Value = 212.3149906599047
$ws.Range("R5").Value = 65.71641987298243
$ws.Range("I6").Value = 10.12574714858493
$ws.Range("J6").Value = 126.0910353404088
$ws.Range("R6").Value = 45.52166981132082
$ws.Range("J7").Value = 33.63624132272333
$ws.Range("K7").Value = 106.7437663446525
$ws.Range("P7").Value = 135.0065633140411
$ws.Range("Q7").Value = 65.34295837775146
$ws.Range("J8").Value = 169.0966151720738
$ws.Range("Q8").Value = 212.3149906599047
$ws.Range("R8").Value = 65.71641987298243
$ws.Range("I9").Value = 10.12574714858493
$ws.Range("J9").Value = 126.0910353404088
$ws.Range("M9").Value = 142.1340339220183
$ws.Range("R9").Value = 45.52166981132082
$ws.Range("J10").Value = 33.63624132272333
$ws.Range("K10").Value = 106.7437663446525
$ws.Range("P10").Value = 135.0065633140411
$ws.Range("Q10").Value = 65.34295837775146
$ws.Range("J11").Value = 169.0966151720738
$ws.Range("Q11").Value = 212.3149906599047
$ws.Range("R11").Value = 65.71641987298243
$ws.Range("I12").Value = 10.12574714858493
$ws.Range("J12").Value = 126.0910353404088
$ws.Range("K12").Value = 137.841438974359
$ws.Range("P12").Value = 133.9744074143302
$ws.Range("Q12").Value = 139.9817740860215
$ws.Range("R12").Value = 45.52166981132082
$ws.Range("K13").Value = 106.7437663446525
$ws.Range("P13").Value = 135.0065633140411
$ws.Range("Q13").Value = 65.34295837775146
$ws.Range("J14").Value = 169.0966151720738
$ws.Range("Q14").Value = 212.3149906599047
$ws.Range("I15").Value = 10.12574714858493
$ws.Range("J15").Value = 126.0910353404088
$ws.Range("R15").Value = 45.52166981132082
$ws.Range("J16").Value = 33.63624132272333
$ws.Range("K16").Value = 106.7437663446525
$ws.Range("P16").Value = 135.0065633140411
$ws.Range("Q16").Value = 65.34295837775146
$ws.Range("J17").Value = 169.0966151720738
$ws.Range("Q17").Value = 212.3149906599047
$ws.Range("R17").Value = 65.71641987298243
$ws.Range("I18").Value = 10.12574714858493
$ws.Range("J18").Value = 126.0910353404088
$ws.Range("P18").Value = 133.9744074143302
$ws.Range("R18").Value = 45.52166981132082
$ws.Range("J19").Value = 33.63624132272333
$ws.Range("K19").Value = 106.7437663446525
$ws.Range("L19").Value = 134.8846762812383
$ws.Range("N19").Value = 127.6855444652332
$ws.Range("O19").Value = 138.4565384518428
$ws.Range("P19").Value = 135.0065633140411
$ws.Range("Q19").Value = 65.34295837775146
$ws.Range("I21").Value = 10.12574714858493
$ws.Range("J21").Value = 126.0910353404088
$ws.Range("R21").Value = 45.52166981132082
$ws.Range("J22").Value = 33.63624132272333
$ws.Range("K22").Value = 106.7437663446525
$ws.Range("L22").Value = 134.8846762812383
$ws.Range("M22").Value = 138.9257839476051
$ws.Range("N22").Value = 127.6855444652332
$ws.Range("P22").Value = 135.0065633140411
$ws.Range("Q22").Value = 65.34295837775146
$ws.Range("J23").Value = 169.0966151720738
$ws.Range("Q23").Value = 212.3149906599047
$ws.Range("R23").Value = 65.71641987298243
$ws.Range("I24").Value = 10.12574714858493
$ws.Range("J24").Value = 126.0910353404088
$ws.Range("K24").Value = 137.841438974359
$ws.Range("R24").Value = 45.52166981132082
$ws.Range("J25").Value = 33.63624132272333
$ws.Range("K25").Value = 106.7437663446525
$ws.Range("L25").Value = 134.8846762812383
$ws.Range("M25").Value = 138.9257839476051
$ws.Range("Q25").Value = 65.34295837775146
$ws.Range("J26").Value = 169.0966151720738
$ws.Range("O26").Value = 230.0982114216867
$ws.Range("Q26").Value = 212.3149906599047
$ws.Range("R26").Value = 65.71641987298243
$ws.Range("I27").Value = 10.12574714858493
$ws.Range("J27").Value = 126.0910353404088
$ws.Range("R27").Value = 45.52166981132082
$ws.Range("J28").Value = 33.63624132272333
$ws.Range("K28").Value = 106.7437663446525
$ws.Range("P28").Value = 135.0065633140411
$ws.Range("Q28").Value = 65.34295837775146
$ws.Range("J29").Value = 169.0966151720738
$ws.Range("Q29").Value = 212.3149906599047
$ws.Range("R29").Value = 65.71641987298243
$ws.Range("I30").Value = 10.12574714858493
$ws.Range("J30").Value = 126.0910353404088
$ws.Range("K30").Value = 137.841438974359
$ws.Range("P30").Value = 133.9744074143302
$ws.Range("Q30").Value = 139.9817740860215
$ws.Range("R30").Value = 45.52166981132082
$ws.Range("J31").Value = 33.63624132272333
$ws.Range("P31").Value = 135.0065633140411
$ws.Range("Q31").Value = 65.34295837775146
$ws.Range("J32").Value = 169.0966151720738
$ws.Range("Q32").Value = 212.3149906599047
$ws.Range("R32").Value = 65.71641987298243
$ws.Range("I33").Value = 10.12574714858493
$ws.Range("J33").Value = 126.0910353404088
$ws.Range("K33").Value = 137.841438974359
$ws.Range("L33").Value = 138.5543797798742
$ws.Range("P33").Value = 133.9744074143302
$ws.Range("Q33").Value = 139.9817740860215
$ws.Range("R33").Value = 45.52166981132082
$ws.Range("J34").Value = 33.63624132272333
$ws.Range("K34").Value = 106.7437663446525
$ws.Range("P34").Value = 135.0065633140411
$ws.Range("Q34").Value = 65.34295837775146
$ws.Range("J35").Value = 169.0966151720738
$ws.Range("P35").Value = 231.2329957552695
$ws.Range("I36").Value = 10.12574714858493
$ws.Range("J36").Value = 126.0910353404088
$ws.Range("K36").Value = 137.841438974359
$ws.Range("R36").Value = 45.52166981132082
$ws.Range("J37").Value = 33.63624132272333
$ws.Range("K37").Value = 106.7437663446525
$ws.Range("L37").Value = 134.8846762812383
$ws.Range("M37").Value = 138.9257839476051
$ws.Range("N37").Value = 127.6855444652332
$ws.Range("O37").Value = 138.4565384518428
$ws.Range("P37").Value = 135.0065633140411
$ws.Range("Q37").Value = 65.34295837775146
$ws.Range("J38").Value = 169.0966151720738
$ws.Range("K38").Value = 220.0898510449805
$ws.Range("Q38").Value = 212.3149906599047
$ws.Range("R38").Value = 65.71641987298243
$ws.Range("I39").Value = 10.12574714858493
$ws.Range("J39").Value = 126.0910353404088
$ws.Range("K39").Value = 137.841438974359
$ws.Range("O39").Value = 142.5962444444444
$ws.Range("P39").Value = 133.9744074143302
$ws.Range("Q39").Value = 139.9817740860215
$ws.Range("R39").Value = 45.52166981132082
$ws.Range("P40").Value = 135.0065633140411
$ws.Range("Q40").Value = 65.34295837775146
$ws.Range("J41").Value = 169.0966151720738
$ws.Range("L41").Value = 235.7664149699872
$ws.Range("O41").Value = 230.0982114216867
$ws.Range("Q41").Value = 212.3149906599047
$ws.Range("R41").Value = 65.71641987298243
$ws.Range("I42").Value = 10.12574714858493
$ws.Range("J42").Value = 126.0910353404088
$ws.Range("K42").Value = 137.841438974359
$ws.Range("P42").Value = 133.9744074143302
$ws.Range("R42").Value = 45.52166981132082
$ws.Range("J43").Value = 33.63624132272333
$ws.Range("K43").Value = 106.7437663446525
$ws.Range("P43").Value = 135.0065633140411
$ws.Range("Q43").Value = 65.34295837775146
$ws.Range("J44").Value = 169.0966151720738
$ws.Range("R44").Value = 65.71641987298243
$ws.Range("I45").Value = 10.12574714858493
$ws.Range("J45").Value = 126.0910353404088
$ws.Range("K45").Value = 137.841438974359
$ws.Range("J46").Value = 33.63624132272333
$ws.Range("K46").Value = 106.7437663446525
$ws.Range("P46").Value = 135.0065633140411
$ws.Range("Q46").Value = 65.34295837775146

# --- Unmet Demand sheet ---
$ws = $wb.Worksheets.Item("Unmet Demand")
$ws.Range("J2").Value = 11.94928935461252
$ws.Range("R2").Value = 149.8691179411497
$ws.Range("I3").Value = 89.39663285141508
$ws.Range("J3").Value = 0.7465913262578567
$ws.Range("K3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 100.1578341526431
$ws.Range("J4").Value = 93.35918011667277
$ws.Range("K4").Value = 22.26949182588285
$ws.Range("P4").Value = 2.721440735106512
$ws.Range("Q4").Value = 86.16204325169439
$ws.Range("J5").Value = 11.94928935461252
$ws.Range("L5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 9.990699214544804
$ws.Range("R5").Value = 149.8691179411497
$ws.Range("I6").Value = 89.39663285141508
$ws.Range("J6").Value = 0.7465913262578567
$ws.Range("R6").Value = 100.1578341526431
$ws.Range("J7").Value = 93.35918011667277
$ws.Range("K7").Value = 22.26949182588285
$ws.Range("P7").Value = 2.721440735106512
$ws.Range("Q7").Value = 86.16204325169439
$ws.Range("J8").Value = 11.94928935461252
$ws.Range("Q8").Value = 9.990699214544804
$ws.Range("R8").Value = 149.8691179411497
$ws.Range("I9").Value = 89.39663285141508
$ws.Range("J9").Value = 0.7465913262578567
$ws.Range("M9").Value = 0
$ws.Range("R9").Value = 100.1578341526431
$ws.Range("J10").Value = 93.35918011667277
$ws.Range("K10").Value = 22.26949182588285
$ws.Range("P10").Value = 2.721440735106512
$ws.Range("Q10").Value = 86.16204325169439
$ws.Range("J11").Value = 11.94928935461252
$ws.Range("Q11").Value = 9.990699214544804
$ws.Range("R11").Value = 149.8691179411497
$ws.Range("I12").Value = 89.39663285141508
$ws.Range("J12").Value = 0.7465913262578567
$ws.Range("K12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 100.1578341526431
$ws.Range("K13").Value = 22.26949182588285
$ws.Range("P13").Value = 2.721440735106512
$ws.Range("Q13").Value = 86.16204325169439
$ws.Range("J14").Value = 11.94928935461252
$ws.Range("Q14").Value = 9.990699214544804
$ws.Range("I15").Value = 89.39663285141508
$ws.Range("J15").Value = 0.7465913262578567
$ws.Range("R15").Value = 100.1578341526431
$ws.Range("J16").Value = 93.35918011667277
$ws.Range("K16").Value = 22.26949182588285
$ws.Range("P16").Value = 2.721440735106512
$ws.Range("Q16").Value = 86.16204325169439
$ws.Range("J17").Value = 11.94928935461252
$ws.Range("Q17").Value = 9.990699214544804
$ws.Range("R17").Value = 149.8691179411497
$ws.Range("I18").Value = 89.39663285141508
$ws.Range("J18").Value = 0.7465913262578567
$ws.Range("P18").Value = 0
$ws.Range("R18").Value = 100.1578341526431
$ws.Range("J19").Value = 93.35918011667277
$ws.Range("K19").Value = 22.26949182588285
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 2.721440735106512
$ws.Range("Q19").Value = 86.16204325169439
$ws.Range("I21").Value = 89.39663285141508
$ws.Range("J21").Value = 0.7465913262578567
$ws.Range("R21").Value = 100.1578341526431
$ws.Range("J22").Value = 93.35918011667277
$ws.Range("K22").Value = 22.26949182588285
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("P22").Value = 2.721440735106512
$ws.Range("Q22").Value = 86.16204325169439
$ws.Range("J23").Value = 11.94928935461252
$ws.Range("Q23").Value = 9.990699214544804
$ws.Range("R23").Value = 149.8691179411497
$ws.Range("I24").Value = 89.39663285141508
$ws.Range("J24").Value = 0.7465913262578567
$ws.Range("K24").Value = 0
$ws.Range("R24").Value = 100.1578341526431
$ws.Range("J25").Value = 93.35918011667277
$ws.Range("K25").Value = 22.26949182588285
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("Q25").Value = 86.16204325169439
$ws.Range("J26").Value = 11.94928935461252
$ws.Range("O26").Value = 0
$ws.Range("Q26").Value = 9.990699214544804
$ws.Range("R26").Value = 149.8691179411497
$ws.Range("I27").Value = 89.39663285141508
$ws.Range("J27").Value = 0.7465913262578567
$ws.Range("R27").Value = 100.1578341526431
$ws.Range("J28").Value = 93.35918011667277
$ws.Range("K28").Value = 22.26949182588285
$ws.Range("P28").Value = 2.721440735106512
$ws.Range("Q28").Value = 86.16204325169439
$ws.Range("J29").Value = 11.94928935461252
$ws.Range("Q29").Value = 9.990699214544804
$ws.Range("R29").Value = 149.8691179411497
$ws.Range("I30").Value = 89.39663285141508
$ws.Range("J30").Value = 0.7465913262578567
$ws.Range("K30").Value = 0
$ws.Range("P30").Value = 0
$ws.Range("Q30").Value = 0
$ws.Range("R30").Value = 100.1578341526431
$ws.Range("J31").Value = 93.35918011667277
$ws.Range("P31").Value = 2.721440735106512
$ws.Range("Q31").Value = 86.16204325169439
$ws.Range("J32").Value = 11.94928935461252
$ws.Range("Q32").Value = 9.990699214544804
$ws.Range("R32").Value = 149.8691179411497
$ws.Range("I33").Value = 89.39663285141508
$ws.Range("J33").Value = 0.7465913262578567
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("Q33").Value = 0
$ws.Range("R33").Value = 100.1578341526431
$ws.Range("J34").Value = 93.35918011667277
$ws.Range("K34").Value = 22.26949182588285
$ws.Range("P34").Value = 2.721440735106512
$ws.Range("Q34").Value = 86.16204325169439
$ws.Range("J35").Value = 11.94928935461252
$ws.Range("P35").Value = 0
$ws.Range("I36").Value = 89.39663285141508
$ws.Range("J36").Value = 0.7465913262578567
$ws.Range("K36").Value = 0
$ws.Range("R36").Value = 100.1578341526431
$ws.Range("J37").Value = 93.35918011667277
$ws.Range("K37").Value = 22.26949182588285
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("O37").Value = 0
$ws.Range("P37").Value = 2.721440735106512
$ws.Range("Q37").Value = 86.16204325169439
$ws.Range("J38").Value = 11.94928935461252
$ws.Range("K38").Value = 0
$ws.Range("Q38").Value = 9.990699214544804
$ws.Range("R38").Value = 149.8691179411497
$ws.Range("I39").Value = 89.39663285141508
$ws.Range("J39").Value = 0.7465913262578567
$ws.Range("K39").Value = 0
$ws.Range("O39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("Q39").Value = 0
$ws.Range("R39").Value = 100.1578341526431
$ws.Range("P40").Value = 2.721440735106512
$ws.Range("Q40").Value = 86.16204325169439
$ws.Range("J41").Value = 11.94928935461252
$ws.Range("L41").Value = 0
$ws.Range("O41").Value = 0
$ws.Range("Q41").Value = 9.990699214544804
$ws.Range("R41").Value = 149.8691179411497
$ws.Range("I42").Value = 89.39663285141508
$ws.Range("J42").Value = 0.7465913262578567
$ws.Range("K42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("R42").Value = 100.1578341526431
$ws.Range("J43").Value = 93.35918011667277
$ws.Range("K43").Value = 22.26949182588285
$ws.Range("P43").Value = 2.721440735106512
$ws.Range("Q43").Value = 86.16204325169439
$ws.Range("J44").Value = 11.94928935461252
$ws.Range("R44").Value = 149.8691179411497
$ws.Range("I45").Value = 89.39663285141508
$ws.Range("J45").Value = 0.7465913262578567
$ws.Range("K45").Value = 0
$ws.Range("J46").Value = 93.35918011667277
$ws.Range("K46").Value = 22.26949182588285
$ws.Range("P46").Value = 2.721440735106512
$ws.Range("Q46").Value = 86.16204325169439

# --- Household Surplus sheet ---
$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B2").Value = 82278.4226072624
$ws.Range("B3").Value = 189860.3556558178
$ws.Range("B4").Value = 102064.4400156472
$ws.Range("B5").Value = 120926.4176592548
$ws.Range("B6").Value = 82082.35965882245
$ws.Range("B7").Value = 118534.0864665268
$ws.Range("B8").Value = 45296.44367078988
$ws.Range("B9").Value = 107654.4253284569
$ws.Range("B10").Value = 123878.0089144462
$ws.Range("B11").Value = 117812.0370933206
$ws.Range("B12").Value = 132786.9241618362
$ws.Range("B13").Value = 118131.2096339729
$ws.Range("B14").Value = 158207.4216240068
$ws.Range("B15").Value = 177646.3062847685
$ws.Range("B16").Value = 68317.43341954847

# --- Costs and Revenues sheet ---
$ws = $wb.Worksheets.Item("Costs and Revenues")
$ws.Range("B2").Value = 2317.702045275
$ws.Range("C2").Value = 5348.179032558249
$ws.Range("D2").Value = 2875.054648328092
$ws.Range("E2").Value = 3406.377962232532
$ws.Range("F2").Value = 2312.179145318943
$ws.Range("G2").Value = 3338.98835116977
$ws.Range("H2").Value = 1275.956159740561
$ws.Range("I2").Value = 3032.519023336817
$ws.Range("J2").Value = 3489.521377871726
$ws.Range("K2").Value = 3318.648932206217
$ws.Range("L2").Value = 3740.476736953133
$ws.Range("M2").Value = 3327.63970799924
$ws.Range("N2").Value = 4456.547088000194
$ws.Range("O2").Value = 5004.121303796296
$ws.Range("P2").Value = 1924.434744212634
$ws.Range("B4").Value = 3476.553067912499
$ws.Range("C4").Value = 8022.268548837372
$ws.Range("D4").Value = 4312.581972492138
$ws.Range("E4").Value = 5109.566943348798
$ws.Range("F4").Value = 3468.268717978414
$ws.Range("G4").Value = 5008.482526754655
$ws.Range("H4").Value = 1913.934239610841
$ws.Range("I4").Value = 4548.778535005225
$ws.Range("J4").Value = 5234.282066807589
$ws.Range("K4").Value = 4977.973398309325
$ws.Range("L4").Value = 5610.715105429699
$ws.Range("M4").Value = 4991.459561998859
$ws.Range("N4").Value = 6684.82063200029
$ws.Range("O4").Value = 7506.181955694445
$ws.Range("P4").Value = 2886.65211631895
$ws.Range("B6").Value = -54153.64424660709
$ws.Range("C6").Value = -54153.64424660708
$ws.Range("D6").Value = -54153.64424660709
$ws.Range("E6").Value = -20526.04424660709
$ws.Range("F6").Value = -20526.04424660709
$ws.Range("G6").Value = -20526.04424660709
$ws.Range("H6").Value = -20526.04424660709
$ws.Range("I6").Value = -20526.04424660708
$ws.Range("J6").Value = -20526.04424660708
$ws.Range("K6").Value = -20526.04424660708
$ws.Range("L6").Value = -20526.04424660708
$ws.Range("M6").Value = -20526.04424660709
$ws.Range("N6").Value = -20526.04424660708
$ws.Range("O6").Value = -20526.04424660709
$ws.Range("P6").Value = -20526.04424660708

